$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the Hours (E) and Activities (F) columns for rows 11-13
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = "Scrum, meeting,  studying"

$ws.Range("E12").Value = 8
$ws.Range("F12").Value = "Scrum, meeting with database, setting up comments"

$ws.Range("E13").Value = 7
$ws.Range("F13").Value = "Meeting with frontend, research, documenting and finishing comment uploading to database"

# Update the active selection / view to match the new state
[void]$ws.Range("F13").Select()
$excel.ActiveWindow.ScrollRow = 5 | Out-Null
